$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert two new data rows (weekly price updates) into the existing
# table, which currently spans rows 2:45 (header in row 1).
#
# A brand-new record (Provincia del Elquí, fecha 44806) is inserted
# right before the current row 33, pushing everything from old row 33
# down to row 34.
# ------------------------------------------------------------------
$ws.Rows(33).Insert()

$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44806
$ws.Range("D33").NumberFormat = $ws.Range("D34").NumberFormat
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112026
$ws.Range("G33").Value = "Haba"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 14000
$ws.Range("M33").Value = 13500
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Provincia del Elquí"
$ws.Range("P33").Value = 540
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"

# ------------------------------------------------------------------
# A second new record (Provincia del Elquí, fecha 44448) is inserted
# right before what is now row 37 (old row 36), pushing the remaining
# rows down by one more.
# ------------------------------------------------------------------
$ws.Rows(37).Insert()

$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").Value = 44448
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 100112026
$ws.Range("G37").Value = "Haba"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 60
$ws.Range("K37").Value = 14000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 14500
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Provincia del Elquí"
$ws.Range("P37").Value = 580
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"

# ------------------------------------------------------------------
# Two rows further down (now rows 43 and 44, previously 41 and 42)
# have their Origen/Precio-$/Kg swapped relative to a straight
# shift-down of the old data.
# ------------------------------------------------------------------
$ws.Range("O43").Value = "Provincia de Diguillín"
$ws.Range("P43").Value = 270

$ws.Range("O44").Value = "Región del Maule"
$ws.Range("P44").Value = 340

# ------------------------------------------------------------------
# Finally, append two brand-new rows (46 and 47) at the bottom of the
# table with additional historical records.
# ------------------------------------------------------------------
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 44484
$ws.Range("D46").NumberFormat = $ws.Range("D45").NumberFormat
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112026
$ws.Range("G46").Value = "Haba"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 30
$ws.Range("K46").Value = 8500
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = 8750
$ws.Range("N46").Value = "$/saco 25 kilos"
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 350
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"

$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C47").Value = "Ñuble"
$ws.Range("D47").Value = 44516
$ws.Range("D47").NumberFormat = $ws.Range("D45").NumberFormat
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = 100112026
$ws.Range("G47").Value = "Haba"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 7000
$ws.Range("L47").Value = 8000
$ws.Range("M47").Value = 7500
$ws.Range("N47").Value = "$/saco 25 kilos"
$ws.Range("O47").Value = "Provincia de Diguillín"
$ws.Range("P47").Value = 300
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"
